$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test")

# Update the "Lead Time for Changes (Median)" row values:
# B3: "* One Month" -> "One Month"
# C3: "* C" -> "C"
$ws.Range("B3").Value = "One Month"
$ws.Range("C3").Value = "C"

# Update the saved selection to C3
$ws.Range("C3").Select()
